$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.186.18'
$ws.Range('E2').Value = '  -0.99%  '

$ws.Range('D3').Value = '2.631.78'
$ws.Range('E3').Value = '  -1.15%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '602.06'
$ws.Range('E5').Value = '  +1.08%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.30'
$ws.Range('E6').Value = '  -0.36%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('E8').Value = '  -0.94%  '

$ws.Range('D9').Value = '2.630.63'
$ws.Range('E9').Value = '  -1.15%  '

$ws.Range('E10').Value = '  -0.54%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.60'
$ws.Range('E11').Value = '  -0.79%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.362'
$ws.Range('E13').Value = '  +1.75%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.14'
$ws.Range('E14').Value = '  -1.84%  '

$ws.Range('D15').Value = '3.099.98'
$ws.Range('E15').Value = '  -1.22%  '

$ws.Range('D16').Value = '63.024.41'
$ws.Range('E16').Value = '  -1.15%  '

$ws.Range('E17').Value = '  -2.15%  '

$ws.Range('D18').Value = '2.623.64'
$ws.Range('E18').Value = '  -1.46%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.30'
$ws.Range('E19').Value = '  -1.15%  '

$ws.Range('E20').Value = '  +2.21%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '339.91'
$ws.Range('E21').Value = '  -1.10%  '

$ws.Range('E22').Value = '  +1.02%  '

$ws.Range('E23').Value = '  -0.10%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.56'
$ws.Range('E24').Value = '  -3.68%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '66.38'
$ws.Range('E25').Value = '  -2.59%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.63'
$ws.Range('E26').Value = '  -3.24%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.65'
$ws.Range('E27').Value = '  +1.30%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '547.81'
$ws.Range('E28').Value = '  -5.45%  '

$ws.Range('E29').Value = '  -7.57%  '

$ws.Range('E30').Value = '  -2.80%  '

$ws.Range('E31').Value = '  -0.04%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.85'
$ws.Range('E32').Value = '  -2.65%  '

$ws.Range('E33').Value = '  -0.36%  '

$ws.Range('E34').Value = '  -3.21%  '

$ws.Range('D35').Value = '0.0₃0803'
$ws.Range('E35').Value = '  -1.96%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.27'
$ws.Range('E36').Value = '  +10.07%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '165.89'
$ws.Range('E37').Value = '  -5.45%  '

$ws.Range('E38').Value = '  -0.05%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.403'
$ws.Range('E39').Value = '  +0.04%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.96'
$ws.Range('E40').Value = '  -1.33%  '

$ws.Range('E41').Value = '  +5.93%  '

$ws.Range('E42').Value = '  +0.04%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '167.90'
$ws.Range('E43').Value = '  -2.05%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.73'
$ws.Range('E44').Value = '  -1.32%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '22.34'
$ws.Range('E45').Value = '  +2.01%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0570'
$ws.Range('E46').Value = '  +3.09%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.622'
$ws.Range('E47').Value = '  -1.80%  '

$ws.Range('E48').Value = '  +0.58%  '

$ws.Range('E49').Value = '  -0.58%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.58'
$ws.Range('E50').Value = '  -0.82%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.76'
$ws.Range('E51').Value = '  +0.91%  '
